# Insert a new weekly price record at row 192 for "Hortaliza, Terminal La
# Palmera de La Serena - Cebollín". Existing rows 192:226 shift down to
# 193:227 (handled automatically by the row Insert), and the freed-up
# row 192 is populated with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 192:226 down to 193:227, leaving a blank row 192 behind
# (inherits formatting from the surrounding rows, e.g. the date style on D).
$ws.Rows.Item(192).Insert()

# Populate the newly inserted row 192 with the new price observation.
$ws.Range("A192").Value = 8
$ws.Range("B192").Value = "Terminal La Palmera de La Serena"
$ws.Range("C192").Value = "Coquimbo"
$ws.Range("D192").Value = 44785
$ws.Range("E192").Value = 4
$ws.Range("F192").Value = 100112037
$ws.Range("G192").Value = "Cebollín"
$ws.Range("H192").Value = "Sin especificar"
$ws.Range("I192").Value = "Primera"
$ws.Range("J192").Value = 2200
$ws.Range("K192").Value = 1400
$ws.Range("L192").Value = 1600
$ws.Range("M192").Value = 1500
$ws.Range("N192").Value = "$/paquete 6 unidades"
$ws.Range("O192").Value = "Provincia del Elquí"
$ws.Range("P192").Value = 250
$ws.Range("Q192").Value = 6
$ws.Range("R192").Value = "Hortaliza"
